$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Replace the old placeholder text (mp_5306_502_topic_4) with the new one
# (MP_5306_502_3). The trailing-space run that followed the old text is
# absorbed/merged by the Find.Execute replace.
$p.Range.Find.Execute("**ID__AFFARS_mp_5306_502_topic_4__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5306_502_3__ID**", 2)

# Drop the now-orphaned trailing space that used to separate the two runs.
$p.Range.Find.Execute("**ID__AFFARS_MP_5306_502_3__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5306_502_3__ID**", 2)

# Indent moves from 120 twips (6pt) to 225 twips (11.25pt).
$p.Format.LeftIndent = 11.25

# Add a paragraph border with 5pt space-from-text on every edge.
$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
